$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "N° of ..." with "Number of ..." across the relevant cells, and
# normalize the two JSON-in-cell values that were previously split across
# multiple rich-text runs into single plain-text strings.

$ws.Range("B2").Value = "Number of Cases"

$ws.Range("H2").Value = '[{"alias":"Confirmed cases", "variable":"Number of Cases", "modifiers":[{"variable":"Case Status", "value":"Confirmed"}]}, {"alias":"Cases at onset of symptomps date", "variable":"Number of Cases", "modifiers":[{"variable":"period type", "value":"Onset of symptomps date"}]}]'

$ws.Range("H3").Value = '[{"15 days Incidence rate":{"Variable":"Incidence Rate", "Calculation Period":"15 days"}}]'

$ws.Range("B5").Value = "Number of PPE Protective equipment"

$ws.Range("B8").Value = "Number of deaths"

$ws.Range("B10").Value = "Number of infected patients"

$ws.Range("B11").Value = "Number of non-infected patients "

$ws.Range("B13").Value = "Number of hospital staff"

$ws.Range("B15").Value = "Number of hospital resources dependencies"

$ws.Range("B20").Value = "Number of test staff"

$ws.Range("B22").Value = "Number of test resources dependencies"

$ws.Range("B34").Value = "Number of index cases studied "

$ws.Range("B36").Value = "Number of clusters found (and cluster type - definition)"

$ws.Range("B54").Value = "Number of calls from people declared as confirm case"

$ws.Range("B69").Value = "Number of people entering to the country (by origin)"

$ws.Range("C4").Select()
